$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Refatoração manual" column header to "Refatoração Guru"
$ws.Range("D1").Value = "Refatoração Guru"

# Fill in the "Guru" analysis columns (C = resultado, D = comparação) for the
# rows that were previously left blank.
$ws.Range("C9").Value  = "Aplicou corretamente"
$ws.Range("D9").Value  = "Igual a referencia"

$ws.Range("C10").Value = "Aplicou corretamente"
$ws.Range("D10").Value = "Diferente da referencia"

$ws.Range("C11").Value = "Aplicou corretamente"
$ws.Range("D11").Value = "Igual a referencia"

$ws.Range("C12").Value = "Aplicou corretamente"
$ws.Range("D12").Value = "Igual a referencia"

$ws.Range("C13").Value = "Aplicou corretamente"
$ws.Range("D13").Value = "Diferente da referencia"

$ws.Range("C14").Value = "Aplicou corretamente"
$ws.Range("D14").Value = "Igual a referencia"

$ws.Range("C15").Value = "Aplicou corretamente"
$ws.Range("D15").Value = "Diferente s referencia"

$ws.Range("C18").Value = "Aplicou corretamente"
$ws.Range("D18").Value = "Diferente da referencia"

$ws.Range("C20").Value = "Aplicou corretamente"
$ws.Range("D20").Value = "Diferente da referência"

$ws.Range("C21").Value = "Aplicou corretamente"
$ws.Range("D21").Value = "Diferente da referência"

$ws.Range("C22").Value = "Aplicou corretamente"
$ws.Range("D22").Value = "Igual a referencia"

$ws.Range("C23").Value = "Não aplicou"
$ws.Range("D23").Value = "-"

$ws.Range("C24").Value = "Aplicou corretamente"
$ws.Range("D24").Value = "Igual a referencia"

$ws.Range("C25").Value = "Aplicou corretamente"
$ws.Range("D25").Value = "Diferente da referência"

$ws.Range("C26").Value = "Aplicou corretamente"
$ws.Range("D26").Value = "Diferente da referência"

# Row-height tweaks that follow the new (longer/shorter) wrapped text
$ws.Rows.Item(15).RowHeight = 21.75
$ws.Rows.Item(19).RowHeight = 18.75
